$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (table "Overview") -- columns A:G
#   A File Name | B Path And Name | C Extension | D Publish URL
#   E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)

# refresh row 2 (file cb648a3a... -> 4c94b8e1...) and regenerate its hyperlink
$wsOv.Range("A2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Range("B2").Value = "e2e\4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("E2").Value = "Handed back: in sync with en-US"
$wsOv.Range("F2").Value = "Handed back: in sync with en-US"
$wsOv.Range("G2").Value = "2016-08-19 17:05:10"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "e2e\4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null

# append row 3 for the new handback file b630bac9...
$loOv.ListRows.Add() | Out-Null
$wsOv.Range("A3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsOv.Range("B3").Value = "e2e\b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-19 17:05:10"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "e2e\b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn") -- columns A:P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("I2").Hyperlinks.Delete()

$wsZh.Range("A2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-19 17:04:58"
$wsZh.Range("I2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsZh.Range("J2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 17:05:30"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/507f3ed2b033526b3be90fafbab6916fad565a1d/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null

$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-19 17:04:58"
$wsZh.Range("I3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsZh.Range("J3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 17:05:30"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/507f3ed2b033526b3be90fafbab6916fad565a1d/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (table "de-de") -- columns A:P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("I2").Hyperlinks.Delete()

$wsDe.Range("A2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-19 17:05:10"
$wsDe.Range("I2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsDe.Range("J2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 17:05:38"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f71eb0125243445e2fef6e1d60a72ac4c35f559/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null

$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-19 17:05:10"
$wsDe.Range("I3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsDe.Range("J3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 17:05:38"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f71eb0125243445e2fef6e1d60a72ac4c35f559/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null

Write-Host "Done applying handback report updates."
